$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Rename "user10@example.com" -> "user11@example.com" everywhere it appears as
# a cell value (both sheets reference the same shared string).
$ws1.Range("A2").Value2 = "user11@example.com"
$ws2.Range("A2").Value2 = "user11@example.com"

# Keep the hyperlinks' display text / tooltip in sync with the new address text.
foreach ($h in $ws1.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = "user11@example.com"
        $h.ScreenTip = "mailto:user11@example.com"
    }
}
foreach ($h in $ws2.Hyperlinks) {
    if ($h.Range.Address() -eq '$A$2') {
        $h.TextToDisplay = "user11@example.com"
        $h.ScreenTip = "mailto:user11@example.com"
    }
}

# Update the selection kept on the (no-longer-active) signup_data sheet.
$ws1.Range("A2").Select()

# Make signin_data the active sheet, with B10 selected.
$ws2.Activate()
$ws2.Range("B10").Select()
